# Commit: "Latest changes pushed on 08-02-2026"
#
# 1. AdminUsers sheet (sheet2.xml): add a "Type" column (C) with header
#    "Type" and the admin row's type "Admin"; also the sample username/
#    password in row 2 changed from "user27" -> "user28".
# 2. A new worksheet "ManageNews" is appended at the end of the workbook
#    (becomes the active/selected sheet) with a small "Message"/"Welcome"
#    table.

$wb = $excel.ActiveWorkbook

# --- AdminUsers: add Type column, bump sample user27 -> user28 ---------
$admin = $wb.Worksheets.Item("AdminUsers")
$admin.Range("C1").Value = "Type"
$admin.Range("C2").Value = "Admin"
$admin.Range("A2").Value = "user28"
$admin.Range("B2").Value = "user28"
$admin.Range("G17").Select()

# --- New ManageNews sheet, appended after the last existing sheet ------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newsSheet = $wb.Worksheets.Add($null, $lastSheet)
$newsSheet.Name = "ManageNews"
$newsSheet.Range("A1").Value = "Message"
$newsSheet.Range("A2").Value = "Welcome"
$newsSheet.Range("B5").Select()
